$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.038481317993446
$ws.Range("D2").Value = 1.042601762076446
$ws.Range("E2").Value = 0.992614727750844
$ws.Range("F2").Value = 1.048153513249709
$ws.Range("I2").Value = 1.042429009588483
$ws.Range("J2").Value = 1.043578521435032
$ws.Range("K2").Value = 1.045377987432357
$ws.Range("L2").Value = 0.9955398523335997
$ws.Range("M2").Value = 1.050914152533261
$ws.Range("N2").Value = 1.018388171980965

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.039528662580185
$ws.Range("D3").Value = 1.043426543600749
$ws.Range("E3").Value = 0.9936372048519299
$ws.Range("F3").Value = 1.049441225103684
$ws.Range("I3").Value = 1.042782960585228
$ws.Range("J3").Value = 1.044270186223669
$ws.Range("K3").Value = 1.046013597376794
$ws.Range("L3").Value = 0.9963617723202687
$ws.Range("M3").Value = 1.052012621155029
$ws.Range("N3").Value = 1.018623697691576

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.040205904737306
$ws.Range("D4").Value = 1.043959768398149
$ws.Range("E4").Value = 0.9942998659930998
$ws.Range("F4").Value = 1.050274398415984
$ws.Range("I4").Value = 1.04301042210973
$ws.Range("J4").Value = 1.044716706344207
$ws.Range("K4").Value = 1.046423768925614
$ws.Range("L4").Value = 0.9968940712668347
$ws.Range("M4").Value = 1.052722792628003
$ws.Range("N4").Value = 1.018775610470451

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.040490508237326
$ws.Range("D5").Value = 1.044183825209557
$ws.Range("E5").Value = 0.994578699834602
$ws.Range("F5").Value = 1.050624651499011
$ws.Range("I5").Value = 1.043105671754552
$ws.Range("J5").Value = 1.044904176559919
$ws.Range("K5").Value = 1.04659593953739
$ws.Range("L5").Value = 0.9971179600053012
$ws.Range("M5").Value = 1.053021204025242
$ws.Range("N5").Value = 1.018839357864108

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.040538288047905
$ws.Range("D6").Value = 1.044221438836655
$ws.Range("E6").Value = 0.994625531979634
$ws.Range("F6").Value = 1.050683459828571
$ws.Range("I6").Value = 1.043121642597193
$ws.Range("J6").Value = 1.044935639174077
$ws.Range("K6").Value = 1.04662483218189
$ws.Range("L6").Value = 0.9971555583673455
$ws.Range("M6").Value = 1.053071300174488
$ws.Range("N6").Value = 1.018850054489723

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.040209708050676
$ws.Range("D7").Value = 1.043962762691737
$ws.Range("E7").Value = 0.994303590798249
$ws.Range("F7").Value = 1.050279078563372
$ws.Range("I7").Value = 1.043011696313967
$ws.Range("J7").Value = 1.044719212299686
$ws.Range("K7").Value = 1.046426070520563
$ws.Range("L7").Value = 0.9968970624462089
$ws.Range("M7").Value = 1.0527267805822
$ws.Range("N7").Value = 1.018776462724803

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.038835368800958
$ws.Range("D8").Value = 1.042880597044506
$ws.Range("E8").Value = 0.9929600610674297
$ws.Range("F8").Value = 1.048588715588769
$ws.Range("I8").Value = 1.042548953926343
$ws.Range("J8").Value = 1.043812487043106
$ws.Range("K8").Value = 1.045593024967338
$ws.Range("L8").Value = 0.9958175282591056
$ws.Range("M8").Value = 1.051285512675739
$ws.Range("N8").Value = 1.018467870132414

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.036410041299786
$ws.Range("D9").Value = 1.040970119960046
$ws.Range("E9").Value = 0.9906006454969559
$ws.Range("F9").Value = 1.045609514634197
$ws.Range("I9").Value = 1.04172151621527
$ws.Range("J9").Value = 1.042206784065362
$ws.Range("K9").Value = 1.044116569354942
$ws.Range("L9").Value = 0.9939188001724441
$ws.Range("M9").Value = 1.048741060329378
$ws.Range("N9").Value = 1.017920346917302

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.034790693592669
$ws.Range("D10").Value = 1.03969405698252
$ws.Range("E10").Value = 0.989033133672735
$ws.Range("F10").Value = 1.043622875209442
$ws.Range("I10").Value = 1.041161783697248
$ws.Range("J10").Value = 1.041130943164717
$ws.Range("K10").Value = 1.043126508839569
$ws.Range("L10").Value = 0.9926553831429383
$ws.Range("M10").Value = 1.047041457424071
$ws.Range("N10").Value = 1.017552806040918

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.034088898660644
$ws.Range("D11").Value = 1.039140931085017
$ws.Range("E11").Value = 0.988355674866747
$ws.Range("F11").Value = 1.042762489837454
$ws.Range("I11").Value = 1.040917485068499
$ws.Range("J11").Value = 1.04066380923252
$ws.Range("K11").Value = 1.042696429608604
$ws.Range("L11").Value = 0.9921088820399291
$ws.Range("M11").Value = 1.046304703437111
$ws.Range("N11").Value = 1.017393055632052

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.033828128191053
$ws.Range("D12").Value = 1.038935387526165
$ws.Range("E12").Value = 0.9881042295826724
$ws.Range("F12").Value = 1.042442878182315
$ws.Range("I12").Value = 1.040826451071263
$ws.Range("J12").Value = 1.040490100575033
$ws.Range("K12").Value = 1.042536471556545
$ws.Range("L12").Value = 0.9919059725120875
$ws.Range("M12").Value = 1.04603091561816
$ws.Range("N12").Value = 1.017333626431016

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.033884068554752
$ws.Range("D13").Value = 1.038979481306365
$ws.Range("E13").Value = 0.9881581567098651
$ws.Range("F13").Value = 1.042511437195164
$ws.Range("I13").Value = 1.040845991331988
$ws.Range("J13").Value = 1.040527370478644
$ws.Range("K13").Value = 1.042570792501322
$ws.Range("L13").Value = 0.9919494934313052
$ws.Range("M13").Value = 1.046089649707536
$ws.Range("N13").Value = 1.017346378308272

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.034067345189299
$ws.Range("D14").Value = 1.039123942587868
$ws.Range("E14").Value = 0.9883348863814464
$ws.Range("F14").Value = 1.042736071188109
$ws.Range("I14").Value = 1.040909966101331
$ws.Range("J14").Value = 1.040649454389624
$ws.Range("K14").Value = 1.042683211662492
$ws.Range("L14").Value = 0.9920921077337197
$ws.Range("M14").Value = 1.046282074604715
$ws.Range("N14").Value = 1.017388145049073

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.034180255661351
$ws.Range("D15").Value = 1.039212938278855
$ws.Range("E15").Value = 0.9884438009545853
$ws.Range("F15").Value = 1.042874472053623
$ws.Range("I15").Value = 1.04094934454358
$ws.Range("J15").Value = 1.040724648542198
$ws.Range("K15").Value = 1.042752449305115
$ws.Range("L15").Value = 0.9921799884222134
$ws.Range("M15").Value = 1.046400617360675
$ws.Range("N15").Value = 1.017413866879746

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.034837256506296
$ws.Range("D16").Value = 1.039730753817491
$ws.Range("E16").Value = 0.9890781214508737
$ws.Range("F16").Value = 1.043679972617449
$ws.Range("I16").Value = 1.041177956268267
$ws.Range("J16").Value = 1.04116191810204
$ws.Range("K16").Value = 1.043155022723159
$ws.Range("L16").Value = 0.9926916645766087
$ws.Range("M16").Value = 1.047090336008503
$ws.Range("N16").Value = 1.017563395430305

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.035249212339091
$ws.Range("D17").Value = 1.040055409598746
$ws.Range("E17").Value = 0.989476357848556
$ws.Range("F17").Value = 1.044185197888482
$ws.Range("I17").Value = 1.041320841068822
$ws.Range("J17").Value = 1.041435860386266
$ws.Range("K17").Value = 1.043407177392224
$ws.Range("L17").Value = 0.9930127773699352
$ws.Range("M17").Value = 1.047522758634005
$ws.Range("N17").Value = 1.017657029199007

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.035489440766645
$ws.Range("D18").Value = 1.040244719682769
$ws.Range("E18").Value = 0.9897087662937556
$ws.Range("F18").Value = 1.044479872228381
$ws.Range("I18").Value = 1.041403997101776
$ws.Range("J18").Value = 1.041595522096675
$ws.Range("K18").Value = 1.043554122178324
$ws.Range("L18").Value = 0.9932001317071769
$ws.Range("M18").Value = 1.047774904974212
$ws.Range("N18").Value = 1.017711586047385

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.035571342583543
$ws.Range("D19").Value = 1.040309260021985
$ws.Range("E19").Value = 0.9897880325774034
$ws.Range("F19").Value = 1.044580346094915
$ws.Range("I19").Value = 1.041432319584429
$ws.Range("J19").Value = 1.041649941575025
$ws.Range("K19").Value = 1.043604204053941
$ws.Range("L19").Value = 0.9932640239640975
$ws.Range("M19").Value = 1.047860867118087
$ws.Range("N19").Value = 1.017730178667263

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.035205019431802
$ws.Range("D20").Value = 1.040020582916262
$ws.Range("E20").Value = 0.9894336180360679
$ws.Range("F20").Value = 1.044130993569418
$ws.Range("I20").Value = 1.041305530150826
$ws.Range("J20").Value = 1.041406481823178
$ws.Range("K20").Value = 1.043380137319945
$ws.Range("L20").Value = 0.9929783193494215
$ws.Range("M20").Value = 1.047476371909778
$ws.Range("N20").Value = 1.01764698919752

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.03401337735582
$ws.Range("D21").Value = 1.039081404773613
$ws.Range("E21").Value = 0.9882828385668249
$ws.Range("F21").Value = 1.042669922820522
$ws.Range("I21").Value = 1.040891135152833
$ws.Range("J21").Value = 1.040613509091625
$ws.Range("K21").Value = 1.042650112754302
$ws.Range("L21").Value = 0.9920501090198102
$ws.Range("M21").Value = 1.046225413725255
$ws.Range("N21").Value = 1.017375848292874

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.033263608306606
$ws.Range("D22").Value = 1.038490396963267
$ws.Range("E22").Value = 0.9875604150241495
$ws.Range("F22").Value = 1.041751136675072
$ws.Range("I22").Value = 1.040628906928739
$ws.Range("J22").Value = 1.040113810985763
$ws.Range("K22").Value = 1.042189916603506
$ws.Range("L22").Value = 0.9914670000341481
$ws.Range("M22").Value = 1.04543816477316
$ws.Range("N22").Value = 1.017204845870823

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.033661126389523
$ws.Range("D23").Value = 1.038803749825246
$ws.Range("E23").Value = 0.9879432794643023
$ws.Range("F23").Value = 1.042238217956436
$ws.Range("I23").Value = 1.040768078651488
$ws.Range("J23").Value = 1.040378817297936
$ws.Range("K23").Value = 1.042433989332185
$ws.Range("L23").Value = 0.991776070289318
$ws.Range("M23").Value = 1.0458555694956
$ws.Range("N23").Value = 1.017295547365035

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.03522498848449
$ws.Range("D24").Value = 1.040036319768212
$ws.Range("E24").Value = 0.9894529299347244
$ws.Range("F24").Value = 1.044155486213725
$ws.Range("I24").Value = 1.041312449070919
$ws.Range("J24").Value = 1.041419757114689
$ws.Range("K24").Value = 1.043392355974657
$ws.Range("L24").Value = 0.9929938892766442
$ws.Range("M24").Value = 1.047497332315609
$ws.Range("N24").Value = 1.017651526021598

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.037037475876338
$ws.Range("D25").Value = 1.041464447933822
$ws.Range("E25").Value = 0.9912096547607049
$ws.Range("F25").Value = 1.046379789377722
$ws.Range("I25").Value = 1.041936856095459
$ws.Range("J25").Value = 1.042622841032748
$ws.Range("K25").Value = 1.044499281232582
$ws.Range("L25").Value = 0.9944092447426414
$ws.Range("M25").Value = 1.049399436547219
$ws.Range("N25").Value = 1.018062339291407
